# Updates the cryptocurrency price/volume table (cols B-E, rows 2-51)
# to reflect the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a temporary "Text" number format on the numeric-looking price cells so that
# Excel stores them as strings (matching the source data) instead of auto-converting
# them to floating point numbers; then restore the default "Normal" style so the
# cell formatting matches the original file.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6306"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07494"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.041"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6823"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001040"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.85"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.313"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "230.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.579"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "159.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.520"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06615"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.458"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.488"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.115"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.117"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.849"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6994"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01872"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.846"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.805"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9367"
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.744"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.121"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1163"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.055"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3958"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05695"
$ws.Range("D51").Style = "Normal"

# Plain text cells (coin names, links, percentages) can be assigned directly.
$ws.Range("D2").Value = "29.616.32"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "1.852.42"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("E10").Value = "  +2.53%  "
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").Value = "1.853.47"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("E17").Value = "  +3.45%  "
$ws.Range("D18").Value = "29.616.26"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("E21").Value = "  +1.96%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("E26").Value = "  -1.82%  "
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("E28").Value = "  +16.26%  "
$ws.Range("E29").Value = "  +3.50%  "
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("E32").Value = "  +1.94%  "
$ws.Range("E33").Value = "  +1.58%  "
$ws.Range("E34").Value = "  -0.66%  "
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("E37").Value = "  +2.21%  "
$ws.Range("D38").Value = "1.264.54"
$ws.Range("E38").Value = "  +1.96%  "
$ws.Range("E39").Value = "  +4.81%  "
$ws.Range("E40").Value = "  +6.62%  "
$ws.Range("E41").Value = "  +3.88%  "
$ws.Range("E42").Value = "  +1.45%  "
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("E45").Value = "  +1.38%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E46").Value = "  +4.24%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("E50").Value = "  -0.84%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("E51").Value = "  -0.10%  "
